$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "70.380.65"
$ws.Cells.Item(2, 5).Value = "  -0.04%  "
$ws.Cells.Item(3, 4).Value = "3.609.69"
$ws.Cells.Item(3, 5).Value = "  -0.54%  "
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "580.53"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.84%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "190.00"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -2.21%  "
$ws.Cells.Item(7, 4).Value = "3.605.45"
$ws.Cells.Item(7, 5).Value = "  -0.45%  "
$ws.Cells.Item(9, 5).Value = "  +0.06%  "
$ws.Cells.Item(10, 5).Value = "  +3.87%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.658"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -1.85%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "56.08"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -3.97%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000312"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +7.28%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "9.72"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -2.52%  "
$ws.Cells.Item(15, 4).Value = "4.189.27"
$ws.Cells.Item(15, 5).Value = "  -0.48%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "19.83"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.23%  "
$ws.Cells.Item(17, 4).Value = "3.608.36"
$ws.Cells.Item(17, 5).Value = "  -0.54%  "
$ws.Cells.Item(18, 4).Value = "70.408.54"
$ws.Cells.Item(18, 5).Value = "  +0.06%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "12.70"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.28%  "
$ws.Cells.Item(20, 5).Value = "  +0.22%  "
$ws.Cells.Item(21, 5).Value = "  -1.85%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "494.39"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.27%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "19.25"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.48%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "4.94"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -6.96%  "
$ws.Cells.Item(25, 5).Value = "  +5.08%  "
$ws.Cells.Item(26, 5).Value = "  -1.83%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "2.99"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -5.53%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "11.11"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -3.67%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "9.48"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -1.22%  "
$ws.Cells.Item(30, 5).Value = "  -2.25%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "7.60"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -3.84%  "
$ws.Cells.Item(32, 5).Value = "  -0.53%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.118"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -2.96%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "65.80"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.18%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "574.44"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -8.35%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "38.58"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -5.55%  "
$ws.Cells.Item(37, 4).Value = "0.0₃0814"
$ws.Cells.Item(37, 5).Value = "  -1.13%  "
$ws.Cells.Item(38, 5).Value = "  +0.17%  "
$ws.Cells.Item(39, 5).Value = "  +16.32%  "
$ws.Cells.Item(40, 5).Value = "  -3.59%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "3.01"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +4.18%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "3.54"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -1.04%  "
$ws.Cells.Item(43, 5).Value = "  -6.22%  "
$ws.Cells.Item(44, 5).Value = "  -4.44%  "
$ws.Cells.Item(45, 4).Value = "3.236.29"
$ws.Cells.Item(45, 5).Value = "  -1.95%  "
$ws.Cells.Item(46, 5).Value = "  -2.15%  "
$ws.Cells.Item(47, 5).Value = "  +6.53%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "3.38"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.92%  "
$ws.Cells.Item(49, 5).Value = "  -0.25%  "
$ws.Cells.Item(50, 5).Value = "  -0.07%  "
$ws.Cells.Item(51, 5).Value = "  -3.82%  "
